$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.072.17"
$ws.Range("E2").Value = "  +3.71%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.267.06"
$ws.Range("E3").Value = "  +1.66%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.28"
$ws.Range("E5").Value = "  -0.65%  "

# Row 6
$ws.Range("E6").Value = "  +0.80%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.11"
$ws.Range("E7").Value = "  +3.50%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.449"
$ws.Range("E9").Value = "  +10.82%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.102"
$ws.Range("E10").Value = "  +12.14%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.02"
$ws.Range("E11").Value = "  -0.94%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.51"
$ws.Range("E12").Value = "  +18.71%  "

# Row 13
$ws.Range("E13").Value = "  +2.20%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.606.96"
$ws.Range("E14").Value = "  +1.87%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.63"
$ws.Range("E15").Value = "  +0.62%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.19"
$ws.Range("E16").Value = "  +9.34%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.842"
$ws.Range("E17").Value = "  +5.39%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.257.12"
$ws.Range("E18").Value = "  +1.50%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.900.69"
$ws.Range("E19").Value = "  +3.66%  "

# Row 20
$ws.Range("E20").Value = "  +7.02%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.55"
$ws.Range("E21").Value = "  +1.67%  "

# Row 22
$ws.Range("E22").Value = "  -1.91%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "252.28"
$ws.Range("E23").Value = "  +3.05%  "

# Row 24
$ws.Range("E24").Value = "  +0.20%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.44"
$ws.Range("E25").Value = "  -2.56%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.29"
$ws.Range("E26").Value = "  -0.33%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.07"
$ws.Range("E27").Value = "  +3.50%  "

# Row 28
$ws.Range("E28").Value = "  +23.99%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.75"
$ws.Range("E29").Value = "  +1.33%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.86"
$ws.Range("E30").Value = "  +2.36%  "

# Row 31
$ws.Range("E31").Value = "  -2.95%  "

# Row 32
$ws.Range("E32").Value = "  -5.30%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.125"
$ws.Range("E33").Value = "  +3.28%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0701"
$ws.Range("E34").Value = "  +7.09%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.79"
$ws.Range("E35").Value = "  +0.70%  "

# Row 36
$ws.Range("E36").Value = "  -1.75%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.84"
$ws.Range("E37").Value = "  +6.87%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.55"
$ws.Range("E38").Value = "  +2.34%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.31"
$ws.Range("E39").Value = "  -2.51%  "

# Row 40
$ws.Range("E40").Value = "  +3.67%  "

# Row 41
$ws.Range("E41").Value = "  -0.13%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.67"
$ws.Range("E42").Value = "  +8.66%  "

# Row 43
$ws.Range("B43").Value = "TerraClassic"
$ws.Range("C43").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.000220"
$ws.Range("E43").Value = "  -2.02%  "

# Row 44
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0973"
$ws.Range("E44").Value = "  +1.08%  "

# Row 45
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.25"
$ws.Range("E45").Value = "  -4.88%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "98.22"
$ws.Range("E46").Value = "  +0.93%  "

# Row 47
$ws.Range("E47").Value = "  -0.76%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.21"
$ws.Range("E48").Value = "  +11.23%  "

# Row 49
$ws.Range("B49").Value = "FTXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.34"
$ws.Range("E49").Value = "  -1.20%  "

# Row 50
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.446.13"
$ws.Range("E50").Value = "  -0.85%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.30"
$ws.Range("E51").Value = "  +3.83%  "
